$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Grade values in column B for specific rows (shared string text values)
$ws.Range("B7").Value = "C"
$ws.Range("B22").Value = "C"
$ws.Range("B39").Value = "C"
$ws.Range("B46").Value = "C"
$ws.Range("B48").Value = "C"
$ws.Range("B50").Value = "C"
$ws.Range("B61").Value = "C"
$ws.Range("B70").Value = "D"

# Update the selected cell/range shown in the sheet view
$ws.Range("B49").Select()
